$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Feathers sheet: append feather-ID rows (P10..P1, S1..S20, alula) for the
# "col_liv" species, mirroring the block structure already used for the
# other species on this sheet.
# ---------------------------------------------------------------------------
$wsFeathers = $wb.Worksheets.Item("Feathers")

$featherIds = @("P10","P9","P8","P7","P6","P5","P4","P3","P2","P1", `
                "S1","S2","S3","S4","S5","S6","S7","S8","S9","S10", `
                "S11","S12","S13","S14","S15","S16","S17","S18","S19","S20", `
                "alula")

$startRow = 33
for ($i = 0; $i -lt $featherIds.Length; $i++) {
    $r = $startRow + $i
    $wsFeathers.Cells.Item($r, 1).Value = "col_liv"
    $wsFeathers.Cells.Item($r, 2).Value = $featherIds[$i]
}

# ---------------------------------------------------------------------------
# Bones sheet: insert a new 9-row block (for "col_liv") above the existing
# "but_jam" block, listing the same bone names used elsewhere on the sheet.
# ---------------------------------------------------------------------------
$wsBones = $wb.Worksheets.Item("Bones")

$boneNames = @("Humerus","Ulna","Radius","Carpometacarpus","Ulnare","Radiale","Digit 2","Digit 3","Digit 4")

$wsBones.Rows("11:19").Insert()

for ($i = 0; $i -lt $boneNames.Length; $i++) {
    $r = 11 + $i
    $wsBones.Cells.Item($r, 1).Value = "col_liv"
    $wsBones.Cells.Item($r, 2).Value = $boneNames[$i]
}

# ---------------------------------------------------------------------------
# Restore selections / active sheet to match the saved view state.
# ---------------------------------------------------------------------------
$wsFeathers.Range("E63").Select()

$wsBones.Activate()
$wsBones.Range("E16").Select()
